# Update Valais COVID-19 daily figures (data refresh / "Add files via upload").
# Only numeric data cells for existing rows (180-209) and the newly
# completed row 210 are touched; everything else (formulas, formatting)
# is left for Excel to recompute naturally on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O ("Nombre de contacts en cours de quarantaine") revisions ---
$ws.Range("O180").Value = 359
$ws.Range("O181").Value = 364
$ws.Range("O182").Value = 371
$ws.Range("O183").Value = 412
$ws.Range("O184").Value = 432
$ws.Range("O185").Value = 458
$ws.Range("O186").Value = 444
$ws.Range("O187").Value = 425
$ws.Range("O188").Value = 356
$ws.Range("O189").Value = 315
$ws.Range("O190").Value = 295
$ws.Range("O191").Value = 258
$ws.Range("O192").Value = 220
$ws.Range("O193").Value = 230
$ws.Range("O194").Value = 239
$ws.Range("O195").Value = 210
$ws.Range("O196").Value = 211
$ws.Range("O197").Value = 222
$ws.Range("O198").Value = 211
$ws.Range("O199").Value = 220

# --- Rows 200-209: later-stage revisions across N/O/P (and a couple of C) ---
$ws.Range("P200").Value = 278

$ws.Range("N201").Value = 78
$ws.Range("P201").Value = 286

$ws.Range("P202").Value = 298

$ws.Range("O203").Value = 286
$ws.Range("P203").Value = 306

$ws.Range("O204").Value = 328
$ws.Range("P204").Value = 332

$ws.Range("N205").Value = 104
$ws.Range("O205").Value = 383
$ws.Range("P205").Value = 340

$ws.Range("C206").Value = 12
$ws.Range("N206").Value = 120
$ws.Range("O206").Value = 403
$ws.Range("P206").Value = 380

$ws.Range("N207").Value = 114
$ws.Range("O207").Value = 421
$ws.Range("P207").Value = 419

$ws.Range("C208").Value = 6
$ws.Range("N208").Value = 110
$ws.Range("O208").Value = 417
$ws.Range("P208").Value = 407

$ws.Range("C209").Value = 7
$ws.Range("N209").Value = 107
$ws.Range("O209").Value = 382
$ws.Range("P209").Value = 402

# --- Row 210 (2020-09-22): newly completed with the day's figures ---
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 10
$ws.Range("I210").Value = 0
$ws.Range("L210").Value = "0"
$ws.Range("M210").Value = "0"
$ws.Range("N210").Value = 98
$ws.Range("O210").Value = 356
$ws.Range("P210").Value = 364

# Restore the on-screen selection to match the author's last position.
$ws.Range("P132").Select()
